$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 271 (shifts existing rows 271-323 down to 272-324)
$ws.Rows("271:271").Insert()

# Populate the new row 271 with the new data record
$ws.Cells.Item(271, 1).Value = 10
$ws.Cells.Item(271, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(271, 3).Value = "La Araucanía"
$ws.Cells.Item(271, 4).Value = 44637
$ws.Cells.Item(271, 5).Value = 9
$ws.Cells.Item(271, 6).Value = 100112040
$ws.Cells.Item(271, 7).Value = "Cilantro"
$ws.Cells.Item(271, 8).Value = "Sin especificar"
$ws.Cells.Item(271, 9).Value = "Primera"
$ws.Cells.Item(271, 10).Value = 100
$ws.Cells.Item(271, 11).Value = 4000
$ws.Cells.Item(271, 12).Value = 4000
$ws.Cells.Item(271, 13).Value = 4000
$ws.Cells.Item(271, 14).Value = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(271, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(271, 16).Value = 2000
$ws.Cells.Item(271, 17).Value = 2
$ws.Cells.Item(271, 18).Value = "Hortaliza"

# Match number format of column D (dates) for the new row
$ws.Cells.Item(271, 4).NumberFormat = $ws.Cells.Item(272, 4).NumberFormat
